# Update cryptos list figures (price + 1h volume change) per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.885.36'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.10%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.586.29'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.18%  '
$ws.Range("E4").Value = '  -0.43%  '
$ws.Range("E5").Value = '  -0.18%  '
$ws.Range("E6").Value = '  -0.41%  '
$ws.Range("E7").Value = '  -2.10%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.248'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.82%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0614'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.30%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.17'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.88%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0789'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.09%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.805.57'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.25%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.577.46'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.98%  '
$ws.Range("E14").Value = '  -1.00%  '
$ws.Range("E15").Value = '  -1.81%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '25.881.97'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.05%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0₃0725'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.48%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '60.17'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.70%  '
$ws.Range("E19").Value = '  -0.39%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '194.47'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.17%  '
$ws.Range("E21").Value = '  +0.03%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.37'
$ws.Range("D22").Style = "Normal"
$ws.Range("E23").Value = '  -0.04%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.130'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.03%  '
$ws.Range("E25").Value = '  -1.14%  '
$ws.Range("E26").Value = '  -0.48%  '
$ws.Range("E27").Value = '  -0.70%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.10'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.64%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.46'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.47%  '
$ws.Range("E30").Value = '  -4.31%  '
$ws.Range("E31").Value = '  +0.23%  '
$ws.Range("E32").Value = '  +1.32%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.02'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.42%  '
$ws.Range("E34").Value = '  +1.99%  '
$ws.Range("E35").Value = '  -2.21%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.095.86'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.10%  '
$ws.Range("E37").Value = '  -0.47%  '
$ws.Range("E38").Value = '  -1.73%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0152'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.03%  '
$ws.Range("E40").Value = '  -0.39%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.781'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.58%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.798'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +6.71%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '93.14'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.32%  '
$ws.Range("E44").Value = '  +0.77%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.718.65'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.27%  '
$ws.Range("E46").Value = '  -1.79%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.52'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.03%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '53.19'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.65%  '
$ws.Range("E49").Value = '  -0.77%  '
$ws.Range("E50").Value = '  -0.84%  '
$ws.Range("E51").Value = '  -0.46%  '
